$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.725.25"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").Value = "3.395.25"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'595.47"
$ws.Range("E5").Value = "  +7.56%  "
$ws.Range("D6").Value = "'186.92"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +4.40%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +4.60%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").Value = "'47.62"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "'0.0000281"
$ws.Range("E12").Value = "  +6.96%  "
$ws.Range("D13").Value = "3.940.48"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "'641.32"
$ws.Range("E14").Value = "  +11.35%  "
$ws.Range("D15").Value = "'8.63"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "68.764.91"
$ws.Range("E16").Value = "  +4.54%  "
$ws.Range("D17").Value = "3.393.16"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("D19").Value = "'18.07"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").Value = "'11.15"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").Value = "'0.915"
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("D24").Value = "'100.24"
$ws.Range("E25").Value = "  +4.46%  "
$ws.Range("E26").Value = "  +6.89%  "
$ws.Range("D27").Value = "'9.81"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").Value = "'32.99"
$ws.Range("E28").Value = "  +8.43%  "
$ws.Range("D29").Value = "'8.73"
$ws.Range("E29").Value = "  +4.18%  "
$ws.Range("D30").Value = "'6.88"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("D31").Value = "'615.50"
$ws.Range("E31").Value = "  +8.88%  "
$ws.Range("D32").Value = "'3.80"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").Value = "4.026.54"
$ws.Range("E33").Value = "  +8.19%  "
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("D37").Value = "'56.83"
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("E38").Value = "  +8.38%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.131"
$ws.Range("E39").Value = "  +4.34%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.31"
$ws.Range("E40").Value = "  +5.86%  "
$ws.Range("D41").Value = "'33.83"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "0.0₃0710"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").Value = "'3.44"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").Value = "'0.345"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("D45").Value = "'0.0424"
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "'2.61"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'1.39"
$ws.Range("E48").Value = "  +13.06%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").Value = "'130.17"
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").Value = "'7.83"
$ws.Range("E51").Value = "  +7.81%  "
